# Apply odds updates + refreshed snapshot timestamps to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTS = "2026-02-22 03:51:59"

# --- Row 2 ---
$ws.Range("F2").Value = 3.1
$ws.Range("I2").Value = 2.32
$ws.Range("J2").Value = 2.6
$ws.Range("BH2").Value = $newTS

# --- Row 3 ---
$ws.Range("Q3").Value = 1.44
$ws.Range("BH3").Value = $newTS

# --- Row 4 ---
$ws.Range("BH4").Value = $newTS

# --- Row 5 ---
$ws.Range("I5").Value = 1.45
$ws.Range("BH5").Value = $newTS

# --- Row 6 ---
$ws.Range("H6").Value = 1.46
$ws.Range("BH6").Value = $newTS

# --- Row 7 ---
$ws.Range("F7").Value = 8.6
$ws.Range("BH7").Value = $newTS

# --- Row 8 ---
$ws.Range("P8").Value = 1.58
$ws.Range("Q8").Value = 2.3
$ws.Range("BH8").Value = $newTS

# --- Row 9 ---
$ws.Range("F9").Value = 3.7
$ws.Range("H9").Value = 2.38
$ws.Range("I9").Value = 2.92
$ws.Range("J9").Value = 2.48
$ws.Range("K9").Value = 3.25
$ws.Range("P9").Value = 1.43
$ws.Range("Q9").Value = 2.98
$ws.Range("BH9").Value = $newTS

# --- Row 10 ---
$ws.Range("F10").Value = 2.3
$ws.Range("I10").Value = 4.7
$ws.Range("K10").Value = 3.3
$ws.Range("BH10").Value = $newTS

# --- Row 11 ---
$ws.Range("BH11").Value = $newTS

# --- Row 12 ---
$ws.Range("BH12").Value = $newTS

# --- Row 13 ---
$ws.Range("G13").Value = 3.9
$ws.Range("P13").Value = 1.44
$ws.Range("Q13").Value = 2.6
$ws.Range("BH13").Value = $newTS
